$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook tracks a rolling weekly window of price observations for
# "Betarraga" at "Vega Monumental Concepción". Each reporting date occupies
# two rows (quality "Primera" then "Segunda"). This week's edit prepends a
# new date's pair of rows at the top of the data block (rows 66-67),
# pushing every existing row down by two positions; the two oldest rows
# that fall off the bottom of the original range (old rows 202-203) are
# appended as new rows at the end (204-205) rather than discarded.

# Insert two blank rows right before the first data row (row 66); Excel
# shifts rows 66:203 down to 68:205 automatically, which reproduces the
# "shift by two" data pattern for every column without touching it by hand.
$ws.Rows("66:67").Insert()

# Populate the two freshly inserted rows with a new date (44469) while
# keeping the same recurring Primera/Segunda attributes the block uses.
$ws.Cells.Item(66, 1).Value2 = 11
$ws.Cells.Item(66, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(66, 3).Value2 = "Bíobío"
$ws.Cells.Item(66, 4).Value2 = 44469
$ws.Cells.Item(66, 5).Value2 = 8
$ws.Cells.Item(66, 6).Value2 = 100114014
$ws.Cells.Item(66, 7).Value2 = "Betarraga"
$ws.Cells.Item(66, 8).Value2 = "Sin especificar"
$ws.Cells.Item(66, 9).Value2 = "Primera"
$ws.Cells.Item(66, 10).Value2 = 800
$ws.Cells.Item(66, 11).Value2 = 600
$ws.Cells.Item(66, 12).Value2 = 700
$ws.Cells.Item(66, 13).Value2 = 650
$ws.Cells.Item(66, 14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(66, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value2 = 130
$ws.Cells.Item(66, 17).Value2 = 5
$ws.Cells.Item(66, 18).Value2 = "Hortaliza"

$ws.Cells.Item(67, 1).Value2 = 11
$ws.Cells.Item(67, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(67, 3).Value2 = "Bíobío"
$ws.Cells.Item(67, 4).Value2 = 44469
$ws.Cells.Item(67, 5).Value2 = 8
$ws.Cells.Item(67, 6).Value2 = 100114014
$ws.Cells.Item(67, 7).Value2 = "Betarraga"
$ws.Cells.Item(67, 8).Value2 = "Sin especificar"
$ws.Cells.Item(67, 9).Value2 = "Segunda"
$ws.Cells.Item(67, 10).Value2 = 400
$ws.Cells.Item(67, 11).Value2 = 500
$ws.Cells.Item(67, 12).Value2 = 500
$ws.Cells.Item(67, 13).Value2 = 500
$ws.Cells.Item(67, 14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(67, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(67, 16).Value2 = 100
$ws.Cells.Item(67, 17).Value2 = 5
$ws.Cells.Item(67, 18).Value2 = "Hortaliza"
